$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 206, shifting existing rows 206..228 down to 207..229.
$ws.Rows.Item(206).Insert()

# Fill the newly inserted row 206 with the new weekly data record.
$ws.Cells.Item(206, 1).Value = 11
$ws.Cells.Item(206, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(206, 3).Value = "Bíobío"
$ws.Cells.Item(206, 4).Value = Get-Date -Year 2022 -Month 2 -Day 25 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(206, 5).Value = 8
$ws.Cells.Item(206, 6).Value = 100112023
$ws.Cells.Item(206, 7).Value = "Brócoli"
$ws.Cells.Item(206, 8).Value = "Sin especificar"
$ws.Cells.Item(206, 9).Value = "Primera"
$ws.Cells.Item(206, 10).Value = 3000
$ws.Cells.Item(206, 11).Value = 700
$ws.Cells.Item(206, 12).Value = 800
$ws.Cells.Item(206, 13).Value = 767
$ws.Cells.Item(206, 14).Value = "`$/unidad"
$ws.Cells.Item(206, 15).Value = "Región del Maule"
$ws.Cells.Item(206, 16).Value = 767
$ws.Cells.Item(206, 17).Value = 1
$ws.Cells.Item(206, 18).Value = "Hortaliza"
